$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (column D) and Volume(1h) (column E) values
# as scraped by the GitHub Actions workflow on 2023-02-15.
$updates = @{
    "D2" = "298.85"
    "E2" = "1.88%"
    "D3" = "42.21"
    "E3" = "4.32%"
    "D4" = "5.017"
    "E4" = "0.14%"
    "D5" = "0.07529"
    "E5" = "2.45%"
    "D6" = "1.600"
    "E6" = "2.15%"
    "D7" = "0.9366"
    "E7" = "1.22%"
    "D8" = "2.389"
    "E8" = "1.59%"
    "D9" = "0.1190"
    "E9" = "2.22%"
    "E10" = "3.53%"
    "D11" = "0.09015"
    "E11" = "2.58%"
    "D12" = "0.04157"
    "E12" = "-4.95%"
    "E13" = "-0.68%"
    "D14" = "0.001280"
    "E14" = "0.77%"
    "D15" = "0.005897"
    "E15" = "0.15%"
    "E16" = "-0.33%"
    "D17" = "4.361"
    "E17" = "1.58%"
    "D18" = "0.3334"
    "E18" = "0.83%"
    "D19" = "8.319"
    "E19" = "6.68%"
    "D20" = "0.1410"
    "E20" = "1.48%"
    "E21" = "12.00%"
    "D22" = "0.04092"
    "E22" = "4.66%"
    "E23" = "0.24%"
    "D24" = "0.003895"
    "E24" = "6.13%"
    "E25" = "8.40%"
    "D38" = "0.02409"
    "E38" = "3.34%"
    "D39" = "0.05225"
    "E39" = "2.47%"
    "D40" = "0.006608"
    "E40" = "19.22%"
    "D41" = "0.007762"
    "E41" = "-1.13%"
    "E42" = "2.59%"
    "D43" = "0.007392"
    "E43" = "0.25%"
    "D44" = "0.007129"
    "E44" = "-11.73%"
    "D45" = "0.2995"
    "E45" = "2.76%"
    "D46" = "0.00006229"
    "E46" = "-0.10%"
    "E47" = "0.02%"
    "D48" = "0.04538"
    "E48" = "-4.45%"
    "E49" = "-0.02%"
    "D50" = "0.00002100"
    "E50" = "0.02%"
    "D51" = "0.0002000"
    "E51" = "0.02%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so values like "298.85" or "1.88%" are not
    # reinterpreted by Excel as numbers/percentages.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
